$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Copy formatting from the last existing data row (row 79) down into the two new rows
$srcRow = $ws.Range("A79:O79")
$srcRow.Copy()
$ws.Range("A80:O80").PasteSpecial(-4122)
$ws.Range("A81:O81").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$rows = @(
    @(45716.900231481479, 8, 6, 229, 393, 369, 388, 2681, 388, 2022, 207, 304, 30, 3189, 4532),
    @(45716.901041666664, 8, 6, 229, 393, 369, 388, 2681, 388, 2022, 207, 304, 30, 3189, 4532)
)

$startRow = 80
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($c = 1; $c -le $data.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $data[$c - 1]
    }
}

$wb.Save()
